# Generate Report for Handoff
# Replaces the two sample files that were previously "handed back" with a new
# pair of files that are now "ready for handoff" (new GUID-based file names,
# new xliff hashes/timestamps), and narrows a few columns that no longer need
# to show long hyperlink text.

$wb = $excel.ActiveWorkbook

$oldFile1 = "383509dd-8e05-4b7e-8273-1d9fa6c3dfe6"
$oldFile2 = "42647d57-8228-4722-a6e3-4fd76a0d03a6"
$newFile1 = "21d41568-f00e-4aaf-90d5-3c89a0e5ceb8"
$newFile2 = "fffff0741068-fb15-4787-8417-99839806d122"

$statusText = "Ready for handoff"
$genDate = "2016-08-24 17:05:31"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newFile1.md"
$wsOverview.Range("B2").Value = "e2e\$newFile1.md"
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("G2").Value = $genDate

$wsOverview.Range("A3").Value = "$newFile2.md"
$wsOverview.Range("B3").Value = "e2e\$newFile2.md"
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText
$wsOverview.Range("G3").Value = $genDate

# Rebuild the hyperlinks so their visible text matches the new file names,
# keeping the same target URLs as before.
$ovLink1Address = $wsOverview.Hyperlinks.Item(1).Address()
$ovLink2Address = $wsOverview.Hyperlinks.Item(2).Address()
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $ovLink1Address, "", "", "e2e\$newFile1.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $ovLink2Address, "", "", "e2e\$newFile2.md")

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhLink1Address = $wsZh.Hyperlinks.Item(1).Address()
$zhLink2Address = $wsZh.Hyperlinks.Item(3).Address()

$wsZh.Range("A2").Value = "$newFile1.md"
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("G2").Value = "$newFile1.885cd247df98f70f27b41408080dead45342786e.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-24 17:05:26"
$wsZh.Range("I2").Value = ""
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"

$wsZh.Range("A3").Value = "$newFile2.md"
$wsZh.Range("C3").Value = $statusText
$wsZh.Range("F3").Value = "True"
$wsZh.Range("G3").Value = "$newFile1.885cd247df98f70f27b41408080dead45342786e.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-24 17:05:26"
$wsZh.Range("I3").Value = ""
$wsZh.Range("I3").Style = "Normal"
$wsZh.Range("J3").Value = ""
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"

# Rebuild hyperlinks: keep A2/A3 (pointing at the same targets as before,
# with refreshed display text); drop the I2/I3 hyperlinks entirely.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zhLink1Address, "", "", "$newFile1.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $zhLink2Address, "", "", "$newFile2.md")

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deLink1Address = $wsDe.Hyperlinks.Item(1).Address()
$deLink2Address = $wsDe.Hyperlinks.Item(3).Address()

$wsDe.Range("A2").Value = "$newFile1.md"
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("G2").Value = "$newFile1.885cd247df98f70f27b41408080dead45342786e.de-de.xlf"
$wsDe.Range("H2").Value = $genDate
$wsDe.Range("I2").Value = ""
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDe.Range("A3").Value = "$newFile2.md"
$wsDe.Range("C3").Value = $statusText
$wsDe.Range("G3").Value = "$newFile1.885cd247df98f70f27b41408080dead45342786e.de-de.xlf"
$wsDe.Range("H3").Value = $genDate
$wsDe.Range("I3").Value = ""
$wsDe.Range("I3").Style = "Normal"
$wsDe.Range("J3").Value = ""
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $deLink1Address, "", "", "$newFile1.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $deLink2Address, "", "", "$newFile2.md")

# ---------------------------------------------------------------------------
# Column width tweaks (columns that used to show long hyperlink text are
# narrowed now that the text they hold is shorter).
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333336   # E
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333336   # F

foreach ($ws in @($wsZh, $wsDe)) {
    $ws.Columns.Item(3).ColumnWidth = 16.333333333333336   # C
    $ws.Columns.Item(9).ColumnWidth = 17.833333333333336   # I
    $ws.Columns.Item(10).ColumnWidth = 20.833333333333336  # J
}
